# Auto-generated Excel COM-interop script
# Updates leve market-price / profit columns (H-N) per scheduled market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 343.625
$ws.Range("I39").Value = 371.35715
$ws.Range("J39").Value = 149.5
$ws.Range("K39").Value = 1114.07145
$ws.Range("L39").Value = 448.5
$ws.Range("M39").Value = -818.0714499999999
$ws.Range("N39").Value = -1040.5

$ws.Range("H40").Value = 1256.3939
$ws.Range("I40").Value = 1248.2273
$ws.Range("J40").Value = 1272.7273
$ws.Range("K40").Value = 1248.2273
$ws.Range("L40").Value = 1272.7273
$ws.Range("M40").Value = -1073.2273
$ws.Range("N40").Value = -1622.7273

$ws.Range("H41").Value = 241.85
$ws.Range("I41").Value = 131.22223
$ws.Range("J41").Value = 332.36365
$ws.Range("K41").Value = 131.22223
$ws.Range("L41").Value = 332.36365
$ws.Range("M41").Value = 308.77777
$ws.Range("N41").Value = -1212.36365

$ws.Range("H61").Value = 434.2857
$ws.Range("I61").Value = 434.2857
$ws.Range("K61").Value = 1302.8571
$ws.Range("M61").Value = -1130.8571

$ws.Range("H64").Value = 2866.3
$ws.Range("I64").Value = 2620
$ws.Range("J64").Value = 2971.8572
$ws.Range("K64").Value = 2620
$ws.Range("L64").Value = 2971.8572
$ws.Range("M64").Value = -2372
$ws.Range("N64").Value = -3467.8572

$ws.Range("H67").Value = 2866.3
$ws.Range("I67").Value = 2620
$ws.Range("J67").Value = 2971.8572
$ws.Range("K67").Value = 2620
$ws.Range("L67").Value = 2971.8572
$ws.Range("M67").Value = -1762
$ws.Range("N67").Value = -4687.8572

$ws.Range("H74").Value = 3714.9429
$ws.Range("I74").Value = 3565.125
$ws.Range("J74").Value = 4041.818
$ws.Range("K74").Value = 3565.125
$ws.Range("L74").Value = 4041.818
$ws.Range("M74").Value = -2629.125
$ws.Range("N74").Value = -5913.818

$ws.Range("H77").Value = 3714.9429
$ws.Range("I77").Value = 3565.125
$ws.Range("J77").Value = 4041.818
$ws.Range("K77").Value = 17825.625
$ws.Range("L77").Value = 20209.09
$ws.Range("M77").Value = -13145.625
$ws.Range("N77").Value = -29569.09

$ws.Range("H137").Value = 1260.9636
$ws.Range("I137").Value = 1008.30304
$ws.Range("J137").Value = 1639.9546
$ws.Range("K137").Value = 3024.90912
$ws.Range("L137").Value = 4919.8638
$ws.Range("M137").Value = -474.9091200000003
$ws.Range("N137").Value = -10019.8638

$ws.Range("H141").Value = 4395.1797
$ws.Range("I141").Value = 2340.739
$ws.Range("J141").Value = 7348.4375
$ws.Range("K141").Value = 7022.217000000001
$ws.Range("L141").Value = 22045.3125
$ws.Range("M141").Value = -1842.217000000001
$ws.Range("N141").Value = -32405.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 12750
$ws.Range("J54").Value = 12750
$ws.Range("L54").Value = 12750
$ws.Range("N54").Value = -14288

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = ""
$ws.Range("N43").Value = 0

$ws.Range("H99").Value = 1449.6086
$ws.Range("I99").Value = 1078.7
$ws.Range("J99").Value = 1734.9231
$ws.Range("K99").Value = 1078.7
$ws.Range("L99").Value = 1734.9231
$ws.Range("M99").Value = 419.3
$ws.Range("N99").Value = -4730.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1407.2593
$ws.Range("I58").Value = 911.5
$ws.Range("J58").Value = 1616
$ws.Range("K58").Value = 911.5
$ws.Range("L58").Value = 1616
$ws.Range("M58").Value = -708.5
$ws.Range("N58").Value = -2022

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = ""
$ws.Range("N97").Value = 0

$ws.Range("H136").Value = 1407.2593
$ws.Range("I136").Value = 911.5
$ws.Range("J136").Value = 1616
$ws.Range("K136").Value = 2734.5
$ws.Range("L136").Value = 4848
$ws.Range("M136").Value = -184.5
$ws.Range("N136").Value = -9948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1179
$ws.Range("I47").Value = 497.5
$ws.Range("J47").Value = 1633.3334
$ws.Range("K47").Value = 1492.5
$ws.Range("L47").Value = 4900.0002
$ws.Range("M47").Value = -1061.5
$ws.Range("N47").Value = -5762.0002

$ws.Range("H50").Value = 402.94116
$ws.Range("I50").Value = 110
$ws.Range("J50").Value = 732.5
$ws.Range("K50").Value = 330
$ws.Range("L50").Value = 2197.5
$ws.Range("M50").Value = 151
$ws.Range("N50").Value = -3159.5

$ws.Range("H53").Value = 402.94116
$ws.Range("I53").Value = 110
$ws.Range("J53").Value = 732.5
$ws.Range("K53").Value = 330
$ws.Range("L53").Value = 2197.5
$ws.Range("M53").Value = 151
$ws.Range("N53").Value = -3159.5

$ws.Range("H88").Value = 4479.8335
$ws.Range("J88").Value = 4479.8335
$ws.Range("L88").Value = 13439.5005
$ws.Range("N88").Value = -14295.5005

$ws.Range("H91").Value = 4479.8335
$ws.Range("J91").Value = 4479.8335
$ws.Range("L91").Value = 13439.5005
$ws.Range("N91").Value = -16403.5005

$ws.Range("H113").Value = 389.59183
$ws.Range("I113").Value = 370.43478
$ws.Range("J113").Value = 406.53845
$ws.Range("K113").Value = 1111.30434
$ws.Range("L113").Value = 1219.61535
$ws.Range("M113").Value = 1058.69566
$ws.Range("N113").Value = -5559.61535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1126.6666
$ws.Range("I9").Value = 690
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 690
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = -466
$ws.Range("N9").Value = -2448

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = 0

$ws.Range("H12").Value = 3833
$ws.Range("I12").Value = 499
$ws.Range("J12").Value = 5500
$ws.Range("K12").Value = 499
$ws.Range("L12").Value = 5500
$ws.Range("M12").Value = -329
$ws.Range("N12").Value = -5840

$ws.Range("H17").Value = 1034.6666
$ws.Range("I17").Value = 827
$ws.Range("J17").Value = 1450
$ws.Range("K17").Value = 827
$ws.Range("L17").Value = 1450
$ws.Range("M17").Value = -657
$ws.Range("N17").Value = -1790

$ws.Range("H21").Value = 76932300
$ws.Range("J21").Value = 76932300
$ws.Range("L21").Value = 76932300
$ws.Range("N21").Value = -76932648

$ws.Range("H24").Value = 2002.5
$ws.Range("I24").Value = 1005
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 1005
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -662
$ws.Range("N24").Value = -3686

$ws.Range("H25").Value = 1201.4
$ws.Range("I25").Value = 2007
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 2007
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -1777
$ws.Range("N25").Value = -1460

$ws.Range("H30").Value = 625
$ws.Range("I30").Value = 625
$ws.Range("K30").Value = 625
$ws.Range("M30").Value = -517

$ws.Range("H46").Value = 1049.75
$ws.Range("I46").Value = 999.2
$ws.Range("J46").Value = 1134
$ws.Range("K46").Value = 999.2
$ws.Range("L46").Value = 1134
$ws.Range("M46").Value = -811.2
$ws.Range("N46").Value = -1510

$ws.Range("H51").Value = 20004
$ws.Range("J51").Value = 20004
$ws.Range("L51").Value = 20004
$ws.Range("N51").Value = -20960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 31156.857
$ws.Range("J93").Value = 31156.857
$ws.Range("L93").Value = 31156.857
$ws.Range("N93").Value = -36148.857

$ws.Range("H136").Value = 10798.167
$ws.Range("I136").Value = 2369
$ws.Range("J136").Value = 22599
$ws.Range("K136").Value = 7107
$ws.Range("L136").Value = 67797
$ws.Range("M136").Value = -4557
$ws.Range("N136").Value = -72897

